# Apply the OOXML diff to the workbook using the Excel COM object model.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear cells that are no longer used in the new layout ---
# F/G/H columns: rows that used to hold "V spur", "V proboscis", "V dip", "V nectar",
# "V dip + V nectar", "h = 3 * V/(pi * r^2)", "nectar height w/ proboscis",
# "new nectar height minus nectar height zero", "new dip height" labels/values/units.
$ws.Range("F5:H20").ClearContents()

# K20:M20 ("volume dipped" / L20 formula / unit "cm3") no longer used
$ws.Range("K20:M20").ClearContents()

# O5:R20 block (r spur quadrado / r proboscis quadrado / nectar height / dip height /
# h+H helper computations previously in columns O, P, R) is entirely removed
$ws.Range("O5:R20").ClearContents()

# --- B3 / K4 header labels keep the same text (SPUR WIDTH MEASUREMENTS / PROBOSCIS WIDTH MEASUREMENTS) ---
$ws.Range("B3").Value = "SPUR WIDTH MEASUREMENTS"
$ws.Range("K4").Value = "PROBOSCIS WIDTH MEASUREMENTS"

# --- New / relocated helper formulas in column F/G (rows 9-17) ---
$ws.Range("F9").Value = "r spur squared"
$ws.Range("G9").Formula = "=C22^2"

$ws.Range("F10").Value = "r proboscis squared"
$ws.Range("G10").Formula = "=L19^2"

$ws.Range("F11").Value = "nectar height"
$ws.Range("G11").Formula = "=C23"

$ws.Range("F14").Value = "h + H"
$ws.Range("G14").Formula = "=G11*(( G9/(G9-G10))  ^ (1/3))"

$ws.Range("F17").Value = "% proboscis (20.5 cm)"

$ws.Range("F16").Value = "% proboscis (30 cm)"
$ws.Range("G16").Formula = "=G14/L17"

# --- K18/L18 and K19/L19 swap meaning: proboscis length 2 (constant 20.5) now sits at row 18,
#     proboscis radius (=L16/2 formula) moves to row 19 ---
$ws.Range("K18").Value = "proboscis length 2"
$ws.Range("L18").Value = 20.5

$ws.Range("K19").Value = "proboscis radius"
$ws.Range("L19").Formula = "=L16/2"

# --- Column width adjustments (K got wider to fit "proboscis length 2" / "% proboscis (30 cm)";
#     O keeps a stored bestFit width even though its cells are now empty) ---
$ws.Columns.Item(11).ColumnWidth = 19.5
$ws.Columns.Item(15).ColumnWidth = 18.8

# --- Selection / active cell to match the saved view state ---
$ws.Range("K7").Select()
